# Refresh cryptos list with latest values from coinranking.com
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '67.665.45'
$ws.Cells.Item(2, 5).Value = '  -2.37%  '
$ws.Cells.Item(3, 4).Value = '2.421.21'
$ws.Cells.Item(4, 5).Value = '  +0.11%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '550.83'
$ws.Cells.Item(5, 5).Value = '  -2.50%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '158.23'
$ws.Cells.Item(6, 5).Value = '  -3.15%  '
$ws.Cells.Item(7, 5).Value = '  +0.10%  '
$ws.Cells.Item(8, 5).Value = '  -2.92%  '
$ws.Cells.Item(9, 4).Value = '2.419.63'
$ws.Cells.Item(9, 5).Value = '  -2.74%  '
$ws.Cells.Item(10, 5).Value = '  -7.98%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.164'
$ws.Cells.Item(11, 5).Value = '  -1.74%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.332'
$ws.Cells.Item(12, 5).Value = '  -5.37%  '
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '4.69'
$ws.Cells.Item(13, 5).Value = '  -3.97%  '
$ws.Cells.Item(14, 4).Value = '2.869.76'
$ws.Cells.Item(14, 5).Value = '  -2.30%  '
$ws.Cells.Item(15, 4).Value = '67.598.52'
$ws.Cells.Item(15, 5).Value = '  -2.47%  '
$ws.Cells.Item(16, 5).Value = '  -6.25%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '22.77'
$ws.Cells.Item(17, 5).Value = '  -5.94%  '
$ws.Cells.Item(18, 4).Value = '2.443.65'
$ws.Cells.Item(18, 5).Value = '  -1.95%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '10.60'
$ws.Cells.Item(19, 5).Value = '  -4.96%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '335.65'
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '6.92'
$ws.Cells.Item(21, 5).Value = '  -5.87%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '3.70'
$ws.Cells.Item(22, 5).Value = '  -3.81%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '0.998'
$ws.Cells.Item(23, 5).Value = '  -0.29%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '1.80'
$ws.Cells.Item(24, 5).Value = '  -5.48%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '65.75'
$ws.Cells.Item(25, 5).Value = '  -5.13%  '
$ws.Cells.Item(26, 4).Value = '2.548.03'
$ws.Cells.Item(26, 5).Value = '  -2.40%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '3.59'
$ws.Cells.Item(27, 5).Value = '  -7.22%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '0.999'
$ws.Cells.Item(28, 5).Value = '  -0.10%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '7.92'
$ws.Cells.Item(29, 5).Value = '  -8.18%  '
$ws.Cells.Item(30, 4).Value = '0.0₃0794'
$ws.Cells.Item(30, 5).Value = '  -8.36%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '6.95'
$ws.Cells.Item(31, 5).Value = '  -8.96%  '
$ws.Cells.Item(32, 5).Value = '  +0.15%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '416.68'
$ws.Cells.Item(33, 5).Value = '  -5.42%  '
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '1.60'
$ws.Cells.Item(34, 5).Value = '  -5.72%  '
$ws.Cells.Item(35, 5).Value = '  -6.48%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '157.54'
$ws.Cells.Item(36, 5).Value = '  +1.60%  '
$ws.Cells.Item(37, 5).Value = '  -0.32%  '
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.999'
$ws.Cells.Item(38, 5).Value = '  -0.21%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.107'
$ws.Cells.Item(39, 5).Value = '  -5.34%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '17.48'
$ws.Cells.Item(40, 5).Value = '  -3.28%  '
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.296'
$ws.Cells.Item(41, 5).Value = '  -5.45%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '4.26'
$ws.Cells.Item(42, 5).Value = '  -6.58%  '
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '1.42'
$ws.Cells.Item(43, 5).Value = '  -9.80%  '
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '1.05'
$ws.Cells.Item(44, 5).Value = '  -1.48%  '
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '131.47'
$ws.Cells.Item(45, 5).Value = '  -5.07%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '1.96'
$ws.Cells.Item(46, 5).Value = '  -8.40%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '3.26'
$ws.Cells.Item(47, 5).Value = '  -4.75%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '0.0707'
$ws.Cells.Item(48, 5).Value = '  -2.62%  '
$ws.Cells.Item(49, 2).Value = 'ARBITRUM'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '0.466'
$ws.Cells.Item(49, 5).Value = '  -8.68%  '
$ws.Cells.Item(50, 2).Value = 'Mantle'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '0.548'
$ws.Cells.Item(50, 5).Value = '  -3.94%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '0.0897'
$ws.Cells.Item(51, 5).Value = '  -2.38%  '
